$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.300.49"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.094.17"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0847"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "2.404.36"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.775"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "2.109.00"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "38.192.72"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.130"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.25%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  +9.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.05%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").Value = "1.546.21"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0906"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "2.291.61"
$ws.Range("E51").Value = "  +2.92%  "
